$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5349937677383423
$ws.Range("B1").Value = 1.342705965042114
$ws.Range("C1").Value = 2.011679172515869
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 2.415190935134888
